$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 13
$ws.Range("A3").Value = 14

$ws.Range("E12").Select()
